$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header cells: new scrape timestamp + updated row count
$ws1.Cells.Item(2, 1).Value = "Última actualización: 17:35:41"
$ws1.Cells.Item(3, 1).Value = "Total filas: 282"

# Re-sort swaps: adjacent rows whose Hora_Scrap/Linea/Minutos (A/C/D) got
# exchanged while Hora_Llegada/Parada (B/E) stayed put. Use the ORIGINAL
# (pre-insert) row numbers - these pairs are all above the first inserted
# row so their numbering is unaffected by the later inserts.
$swapPairs1 = @(
    @(66, 67),
    @(118, 119),
    @(137, 138),
    @(190, 191),
    @(235, 236),
    @(253, 254),
    @(262, 263),
    @(271, 272)
)
foreach ($pair in $swapPairs1) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $a1 = $ws1.Cells.Item($r1, 1).Value
    $c1 = $ws1.Cells.Item($r1, 3).Value
    $d1 = $ws1.Cells.Item($r1, 4).Value
    $a2 = $ws1.Cells.Item($r2, 1).Value
    $c2 = $ws1.Cells.Item($r2, 3).Value
    $d2 = $ws1.Cells.Item($r2, 4).Value

    $ws1.Cells.Item($r1, 1).Value = $a2
    $ws1.Cells.Item($r1, 3).Value = $c2
    $ws1.Cells.Item($r1, 4).Value = $d2
    $ws1.Cells.Item($r2, 1).Value = $a1
    $ws1.Cells.Item($r2, 3).Value = $c1
    $ws1.Cells.Item($r2, 4).Value = $d1
}

# New scrape rows. Each Insert() pushes the target row (and below) down by
# one, so inserting sequentially at the destination row number lands each
# new row exactly where it belongs in the final sheet.
$ws1.Rows.Item(248).Insert()
$ws1.Cells.Item(248, 1).Value = "17:35:41"
$ws1.Cells.Item(248, 2).Value = "17:37"
$ws1.Cells.Item(248, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(248, 4).Value = 2
$ws1.Cells.Item(248, 5).Value = "LP1912"

$ws1.Rows.Item(276).Insert()
$ws1.Cells.Item(276, 1).Value = "17:35:41"
$ws1.Cells.Item(276, 2).Value = "18:37"
$ws1.Cells.Item(276, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(276, 4).Value = 62
$ws1.Cells.Item(276, 5).Value = "LP1912"

$ws1.Rows.Item(283).Insert()
$ws1.Cells.Item(283, 1).Value = "17:35:41"
$ws1.Cells.Item(283, 2).Value = "19:03"
$ws1.Cells.Item(283, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(283, 4).Value = 88
$ws1.Cells.Item(283, 5).Value = "LP1912"

# Remaining rows are brand-new and simply appended after the old last row
# (which is now at row 284).
$ws1.Cells.Item(285, 1).Value = "17:35:41"
$ws1.Cells.Item(285, 2).Value = "19:16"
$ws1.Cells.Item(285, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(285, 4).Value = 101
$ws1.Cells.Item(285, 5).Value = "LP1912"

$ws1.Cells.Item(286, 1).Value = "17:35:41"
$ws1.Cells.Item(286, 2).Value = "19:17"
$ws1.Cells.Item(286, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(286, 4).Value = 102
$ws1.Cells.Item(286, 5).Value = "LP1912"

$ws1.Cells.Item(287, 1).Value = "17:35:41"
$ws1.Cells.Item(287, 2).Value = "19:27"
$ws1.Cells.Item(287, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(287, 4).Value = 112
$ws1.Cells.Item(287, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 17:35:41"
$ws2.Cells.Item(3, 1).Value = "Total filas: 48"

$ws2.Rows.Item(44).Insert()
$ws2.Cells.Item(44, 1).Value = "17:35:41"
$ws2.Cells.Item(44, 2).Value = "17:37"
$ws2.Cells.Item(44, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(44, 4).Value = 2
$ws2.Cells.Item(44, 5).Value = "LP1912"

$ws2.Rows.Item(51).Insert()
$ws2.Cells.Item(51, 1).Value = "17:35:41"
$ws2.Cells.Item(51, 2).Value = "19:03"
$ws2.Cells.Item(51, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(51, 4).Value = 88
$ws2.Cells.Item(51, 5).Value = "LP1912"

$ws2.Cells.Item(53, 1).Value = "17:35:41"
$ws2.Cells.Item(53, 2).Value = "19:27"
$ws2.Cells.Item(53, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(53, 4).Value = 112
$ws2.Cells.Item(53, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 17:35:41"
$ws3.Cells.Item(3, 1).Value = "Total filas: 42"

$ws3.Cells.Item(47, 1).Value = "17:35:41"
$ws3.Cells.Item(47, 2).Value = "19:23"
$ws3.Cells.Item(47, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(47, 4).Value = 108
$ws3.Cells.Item(47, 5).Value = "L6173"
